$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds in row 2 (Banfield vs Tigre) ---
$ws.Cells.Item(2, 7).Value = 2.55
$ws.Cells.Item(2, 9).Value = 3
$ws.Cells.Item(2, 12).Value = 3.75
$ws.Cells.Item(2, 13).Value = 1.11
$ws.Cells.Item(2, 14).Value = 6.5
$ws.Cells.Item(2, 21).Value = 2.1
$ws.Cells.Item(2, 22).Value = 1.67
$ws.Cells.Item(2, 23).Value = 6.5
$ws.Cells.Item(2, 24).Value = 11
$ws.Cells.Item(2, 26).Value = 26
$ws.Cells.Item(2, 37).Value = 29
$ws.Cells.Item(2, 41).Value = 17

# --- Update existing odds in row 7 (Botafogo SP vs Avai) ---
$ws.Cells.Item(7, 13).Value = 1.13
$ws.Cells.Item(7, 14).Value = 6

# --- Insert a new row at position 8 (shifts old rows 8-10 down to 9-11) ---
$ws.Rows(8).Insert()

# --- Populate new row 8: Ceara vs America MG (BRAZIL - SERIE B) ---
$ws.Cells.Item(8, 1).Value = "rNtoHigg"
$ws.Cells.Item(8, 2).Value = "18/11/2024"
$ws.Cells.Item(8, 3).Value = "21:45"
$ws.Cells.Item(8, 4).Value = "BRAZIL - SERIE B"
$ws.Cells.Item(8, 5).Value = "Ceara"
$ws.Cells.Item(8, 6).Value = "America MG"
$ws.Cells.Item(8, 7).Value = 1.42
$ws.Cells.Item(8, 8).Value = 4.33
$ws.Cells.Item(8, 9).Value = 7.5
$ws.Cells.Item(8, 10).Value = 1.95
$ws.Cells.Item(8, 11).Value = 2.3
$ws.Cells.Item(8, 12).Value = 7
$ws.Cells.Item(8, 13).Value = 1.05
$ws.Cells.Item(8, 14).Value = 11
$ws.Cells.Item(8, 15).Value = 1.25
$ws.Cells.Item(8, 16).Value = 3.75
$ws.Cells.Item(8, 17).Value = 1.85
$ws.Cells.Item(8, 18).Value = 2
$ws.Cells.Item(8, 19).Value = 1.36
$ws.Cells.Item(8, 20).Value = 3
$ws.Cells.Item(8, 21).Value = 2
$ws.Cells.Item(8, 22).Value = 1.73
$ws.Cells.Item(8, 23).Value = 6.5
$ws.Cells.Item(8, 24).Value = 6.5
$ws.Cells.Item(8, 25).Value = 8.5
$ws.Cells.Item(8, 26).Value = 9.5
$ws.Cells.Item(8, 27).Value = 12
$ws.Cells.Item(8, 28).Value = 29
$ws.Cells.Item(8, 29).Value = 11
$ws.Cells.Item(8, 30).Value = 8.5
$ws.Cells.Item(8, 31).Value = 21
$ws.Cells.Item(8, 32).Value = 67
$ws.Cells.Item(8, 33).Value = 17
$ws.Cells.Item(8, 34).Value = 41
$ws.Cells.Item(8, 35).Value = 23
$ws.Cells.Item(8, 36).Value = 81
$ws.Cells.Item(8, 37).Value = 51
$ws.Cells.Item(8, 38).Value = 51
$ws.Cells.Item(8, 39).Value = 351
$ws.Cells.Item(8, 40).Value = 3.25
$ws.Cells.Item(8, 41).Value = 7
$ws.Cells.Item(8, 42).Value = 19
$ws.Cells.Item(8, 43).Value = 21
$ws.Cells.Item(8, 44).Value = 41
$ws.Cells.Item(8, 45).Value = 151
$ws.Cells.Item(8, 46).Value = 3
$ws.Cells.Item(8, 47).Value = 9.5
$ws.Cells.Item(8, 48).Value = 67
$ws.Cells.Item(8, 49).Value = 8.5
$ws.Cells.Item(8, 50).Value = 41
$ws.Cells.Item(8, 51).Value = 41
$ws.Cells.Item(8, 52).Value = 151
$ws.Cells.Item(8, 53).Value = 151
$ws.Cells.Item(8, 54).Value = 351
$ws.Cells.Item(8, 55).Value = 81
$ws.Cells.Item(8, 56).Value = 81
